$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "ColWidthCal"
$vals = @(1,2,5,8,9,9.5,10,11,11.625,14,14.5,15,16,20,20.5,21,21.5,23,23.375,23.5)
for ($i=0; $i -lt $vals.Length; $i++) {
    $ws.Columns.Item($i+1).ColumnWidth = $vals[$i]
}
Write-Output "done"
